$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 121.75
$ws.Range("J2").Value = 100
$ws.Range("L2").Value = 100
$ws.Range("N2").Value = -326
$ws.Range("H19").Value = 968.6
$ws.Range("I19").Value = 960.75
$ws.Range("J19").Value = 1000
$ws.Range("K19").Value = 960.75
$ws.Range("L19").Value = 1000
$ws.Range("M19").Value = -785.75
$ws.Range("N19").Value = -1350
$ws.Range("H41").Value = 197.25
$ws.Range("J41").Value = 346.33334
$ws.Range("L41").Value = 346.33334
$ws.Range("N41").Value = -1226.33334
$ws.Range("H55").Value = 197
$ws.Range("I55").Value = 233
$ws.Range("K55").Value = 233
$ws.Range("M55").Value = -19
$ws.Range("H92").Value = 752.35
$ws.Range("I92").Value = 781.5
$ws.Range("J92").Value = 635.75
$ws.Range("K92").Value = 781.5
$ws.Range("L92").Value = 635.75
$ws.Range("M92").Value = 466.5
$ws.Range("N92").Value = -3131.75
$ws.Range("H94").Value = 1513.7142
$ws.Range("I94").Value = 932
$ws.Range("K94").Value = 932
$ws.Range("M94").Value = -481
$ws.Range("H96").Value = 381.14285
$ws.Range("I96").Value = 267.81818
$ws.Range("J96").Value = 796.6667
$ws.Range("K96").Value = 803.45454
$ws.Range("L96").Value = 2390.0001
$ws.Range("M96").Value = 569.54546
$ws.Range("N96").Value = -5136.0001
$ws.Range("H107").Value = 5957.3335
$ws.Range("I107").Value = 5148.8
$ws.Range("K107").Value = 5148.8
$ws.Range("M107").Value = -3228.8
$ws.Range("H132").Value = 2369.2
$ws.Range("I132").Value = 2369.2
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 7107.599999999999
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -4577.599999999999
$ws.Range("N132").ClearContents()
$ws.Range("H137").Value = 1829.3334
$ws.Range("I137").Value = 1829.3334
$ws.Range("K137").Value = 5488.0002
$ws.Range("M137").Value = -2938.0002
$ws.Range("H138").Value = 3527.2173
$ws.Range("I138").Value = 1999.6
$ws.Range("J138").Value = 3951.5557
$ws.Range("K138").Value = 5998.799999999999
$ws.Range("L138").Value = 11854.6671
$ws.Range("M138").Value = -858.7999999999993
$ws.Range("N138").Value = -22134.6671
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10383.368
$ws.Range("I32").Value = 10383.368
$ws.Range("K32").Value = 10383.368
$ws.Range("M32").Value = -10096.368
$ws.Range("H122").Value = 3325
$ws.Range("I122").Value = 3325
$ws.Range("K122").Value = 9975
$ws.Range("M122").Value = -7525
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2297.1333
$ws.Range("I105").Value = 2471.5
$ws.Range("K105").Value = 2471.5
$ws.Range("M105").Value = -724.5
$ws.Range("H107").Value = 1916.5
$ws.Range("I107").Value = 2049.8
$ws.Range("J107").Value = 1250
$ws.Range("K107").Value = 2049.8
$ws.Range("L107").Value = 1250
$ws.Range("M107").Value = -129.8000000000002
$ws.Range("N107").Value = -5090
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 132.55556
$ws.Range("I7").Value = 136.375
$ws.Range("J7").Value = 102
$ws.Range("K7").Value = 136.375
$ws.Range("L7").Value = 102
$ws.Range("M7").Value = -23.375
$ws.Range("N7").Value = -328
$ws.Range("H16").Value = 980.8
$ws.Range("I16").Value = 976
$ws.Range("K16").Value = 976
$ws.Range("M16").Value = -689
$ws.Range("H22").Value = 300
$ws.Range("I22").Value = 300
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 300
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = 50
$ws.Range("N22").ClearContents()
$ws.Range("H31").Value = 3084
$ws.Range("I31").Value = 2855.25
$ws.Range("K31").Value = 2855.25
$ws.Range("M31").Value = -2560.25
$ws.Range("H34").Value = 3084
$ws.Range("I34").Value = 2855.25
$ws.Range("K34").Value = 2855.25
$ws.Range("M34").Value = -2653.25
$ws.Range("H107").Value = 1749.6428
$ws.Range("I107").Value = 2100.2
$ws.Range("J107").Value = 873.25
$ws.Range("K107").Value = 2100.2
$ws.Range("L107").Value = 873.25
$ws.Range("M107").Value = -180.1999999999998
$ws.Range("N107").Value = -4713.25
$ws.Range("H113").Value = 980.8
$ws.Range("I113").Value = 976
$ws.Range("K113").Value = 976
$ws.Range("M113").Value = 1194
$ws.Range("H132").Value = 2562.8462
$ws.Range("I132").Value = 2012.4286
$ws.Range("J132").Value = 3205
$ws.Range("K132").Value = 6037.2858
$ws.Range("L132").Value = 9615
$ws.Range("M132").Value = -3507.2858
$ws.Range("N132").Value = -14675
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 1503003.9
$ws.Range("I8").Value = 1503003.9
$ws.Range("K8").Value = 4509011.699999999
$ws.Range("M8").Value = -4508872.699999999
$ws.Range("H80").Value = 12833.333
$ws.Range("I80").Value = 1500
$ws.Range("J80").Value = 14250
$ws.Range("K80").Value = 4500
$ws.Range("L80").Value = 42750
$ws.Range("M80").Value = -3564
$ws.Range("N80").Value = -44622
$ws.Range("H83").Value = 12833.333
$ws.Range("I83").Value = 1500
$ws.Range("J83").Value = 14250
$ws.Range("K83").Value = 13500
$ws.Range("L83").Value = 128250
$ws.Range("M83").Value = -8820
$ws.Range("N83").Value = -137610
$ws.Range("H113").Value = 1874.375
$ws.Range("J113").Value = 2166.3333
$ws.Range("L113").Value = 6498.999899999999
$ws.Range("N113").Value = -10838.9999
$ws.Range("H121").Value = 1669.2858
$ws.Range("I121").Value = 765
$ws.Range("J121").Value = 2347.5
$ws.Range("K121").Value = 2295
$ws.Range("L121").Value = 7042.5
$ws.Range("M121").Value = -985
$ws.Range("N121").Value = -9662.5
$ws.Range("H132").Value = 1300
$ws.Range("I132").Value = 1300
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 11700
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -9170
$ws.Range("N132").ClearContents()
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1830.1111
$ws.Range("I102").Value = 1830.1111
$ws.Range("K102").Value = 1830.1111
$ws.Range("M102").Value = -208.1111000000001
$ws.Range("H113").Value = 3000
$ws.Range("I113").Value = 3000
$ws.Range("K113").Value = 3000
$ws.Range("M113").Value = -830
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()
$ws.Range("H119").Value = 150000
$ws.Range("J119").Value = 150000
$ws.Range("L119").Value = 150000
$ws.Range("N119").Value = -159676
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2761
$ws.Range("I16").Value = 2761
$ws.Range("K16").Value = 2761
$ws.Range("M16").Value = -2591
$ws.Range("H20").Value = 1643812.5
$ws.Range("I20").Value = 2500750
$ws.Range("J20").Value = 1358166.6
$ws.Range("K20").Value = 2500750
$ws.Range("L20").Value = 1358166.6
$ws.Range("M20").Value = -2500524
$ws.Range("N20").Value = -1358618.6
$ws.Range("H100").Value = 1769.0769
$ws.Range("I100").Value = 1774.8334
$ws.Range("K100").Value = 1774.8334
$ws.Range("M100").Value = -1233.8334
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 40000
$ws.Range("J64").Value = 40000
$ws.Range("L64").Value = 40000
$ws.Range("N64").Value = -40496
$ws.Range("H67").Value = 40000
$ws.Range("J67").Value = 40000
$ws.Range("L67").Value = 40000
$ws.Range("N67").Value = -41716
$ws.Range("H107").Value = 1262
$ws.Range("I107").Value = 1024.25
$ws.Range("K107").Value = 3072.75
$ws.Range("M107").Value = -1152.75
